$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Duplicate row 2's formatting into row 3 (new employee row added under
#    the existing one), then copy its values over so row 3 ends up as a
#    twin of row 2's current data (name column left blank).
$ws.Range("B2:J2").Copy()
$ws.Range("B3:J3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B3").Value = $ws.Range("B2").Value()
$ws.Range("C3").Value = $ws.Range("C2").Value()
$ws.Range("D3").Value = $ws.Range("D2").Value()
$ws.Range("E3").Value = $ws.Range("E2").Value()
$ws.Range("G3").Value = $ws.Range("G2").Value()
$ws.Range("H3").Value = $ws.Range("H2").Value()
$ws.Range("I3").Value = $ws.Range("I2").Value()
$ws.Range("J3").Value = $ws.Range("J2").Value()

# 2. Update the employee name on row 2, and clear the phone-number field on
#    both row 2 and the newly duplicated row 3.
$ws.Range("A2").Value = "Duc Mạnh"
$ws.Range("F2").ClearContents()
$ws.Range("F3").ClearContents()

# 3. Prepare row 4 as a fresh, blank entry row, carrying over the same
#    number-format styling used on rows above (date on B, text on E/F).
$ws.Range("B2").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("E2:F2").Copy()
$ws.Range("E4:F4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4. Move the active selection down to the next entry cell.
$ws.Range("B5").Select()
